$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (co2)
$ws.Range("C2").Value = 540.3259683840465
$ws.Range("D2").Value = 129.5449376653265
$ws.Range("F2").Value = 447
$ws.Range("G2").Value = 493
$ws.Range("H2").Value = 600

# Row 3 (humidity)
$ws.Range("C3").Value = 41.32955035315693
$ws.Range("D3").Value = 4.836281502138328
$ws.Range("F3").Value = 38.15
$ws.Range("G3").Value = 40.59
$ws.Range("H3").Value = 44.24

# Row 4 (pm25)
$ws.Range("C4").Value = 1.429229449054731
$ws.Range("D4").Value = 2.004718137524667
$ws.Range("F4").Value = 0.53
$ws.Range("G4").Value = 1.02
$ws.Range("H4").Value = 1.84

# Row 5 (pressure)
$ws.Range("C5").Value = 323.0843761917369
$ws.Range("D5").Value = 10.66518236359977
$ws.Range("F5").Value = 316.83
$ws.Range("G5").Value = 324.79
$ws.Range("H5").Value = 331.68

# Row 6 (temperature)
$ws.Range("C6").Value = 20.80065859962376
$ws.Range("D6").Value = 2.579163722577207
$ws.Range("F6").Value = 19.43
$ws.Range("G6").Value = 20.89
$ws.Range("H6").Value = 22.27

# Row 7 (rssi)
$ws.Range("C7").Value = -76.29231265034586
$ws.Range("D7").Value = 22.70534258309803
$ws.Range("G7").Value = -74

# Row 8 (snr)
$ws.Range("C8").Value = 7.705900981683519
$ws.Range("D8").Value = 6.833132541115538

# Row 9 (SF)
$ws.Range("C9").Value = 9.319969256228363
$ws.Range("D9").Value = 1.684923486802206

# Row 10 (frequency)
$ws.Range("C10").Value = 867.830103702681
$ws.Range("D10").Value = 0.4614640650157109

# Row 11 (toa)
$ws.Range("C11").Value = 0.5549253774386244
$ws.Range("D11").Value = 0.5885556454556197

# Row 12 (distance)
$ws.Range("C12").Value = 22.73539255218065
$ws.Range("D12").Value = 12.29254234512506

# Row 13 (c_walls)
$ws.Range("C13").Value = 0.6738321139412262
$ws.Range("D13").Value = 0.7504700985557685

# Row 14 (w_walls)
$ws.Range("C14").Value = 1.826334777165204
$ws.Range("D14").Value = 1.664180072947043

# Row 15 (exp_pl)
$ws.Range("C15").Value = 93.69231265034573
$ws.Range("D15").Value = 22.70534258309803
$ws.Range("G15").Value = 91.40000000000001

# Row 16 (n_power)
$ws.Range("C16").Value = -85.56293756301817
$ws.Range("D16").Value = 20.4547474905328
$ws.Range("F16").Value = -102.0738221927363
$ws.Range("G16").Value = -84.26572375596102
$ws.Range("H16").Value = -67.8707776445072

# Row 17 (esp)
$ws.Range("C17").Value = -77.85703658133461
$ws.Range("D17").Value = 25.06701371063542
$ws.Range("F17").Value = -92.71081852649533
$ws.Range("G17").Value = -74.18978441047734
$ws.Range("H17").Value = -57.25410721860875
